$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used elsewhere in the sheet for the Difficulty column
#   Easy   -> RGB(00B050) -> OLE color 5287936
#   Medium -> RGB(FFC000) -> OLE color 49407
$easyColor = 5287936

# --- Row 13: 511. Game Play Analysis I ---
$ws.Range("A13").Value = "511. Game Play Analysis I"
$ws.Range("B13").Value = "Easy"
$ws.Range("B13").Interior.Color = $easyColor
$ws.Range("C13").Value = "Aggregation"
$ws.Range("D13").Value = "Sort df by player_id and event_date with sort_values(), groupby() with player_id then min() on event_date, then reset index. Then rename() the event_date column to first_login inplace, then return result."
$ws.Range("E13").Value = "https://leetcode.com/problems/game-play-analysis-i/solutions/3899132/pandas-very-simple-3-step-approach/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "
$ws.Hyperlinks.Add($ws.Range("E13"), "https://leetcode.com/problems/game-play-analysis-i/solutions/3899132/pandas-very-simple-3-step-approach/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata")
$ws.Range("E13").Style = "Hyperlink"

# --- Row 14: 586. Customer Placing the Largest Number of Orders ---
$ws.Range("A14").Value = "586. Customer Placing the Largest Number of Orders"
$ws.Range("B14").Value = "Easy"
$ws.Range("B14").Interior.Color = $easyColor
$ws.Range("C14").Value = "Aggregation"
$ws.Range("D14").Value = "Group df by customer_number using groupby(), count num of orders for each customer with count() on grouped df and resetting index, fileter where order count is equal to the max order count, return res for customer with largest num orders."
$ws.Range("E14").Value = "https://leetcode.com/problems/customer-placing-the-largest-number-of-orders/solutions/3899213/pandas-2-liner-very-straightforward/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "
$ws.Hyperlinks.Add($ws.Range("E14"), "https://leetcode.com/problems/customer-placing-the-largest-number-of-orders/solutions/3899213/pandas-2-liner-very-straightforward/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata")
$ws.Range("E14").Style = "Hyperlink"

# Grow the table (Table2) to include the two new rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E14"))

# Move the selection to where the user last clicked
$ws.Range("D17").Select() | Out-Null
